# ClausesTemplate_fixed.docx edit script
# Applies:
#  1. "In Word 2019" -> "In Word"
#  2. Remove the stray _GoBack bookmark after the "...agreement that you are
#     working on." paragraph.
#  3. Remove the "You can find an appendix ..." paragraph (duplicate of the
#     "You can find lists ..." paragraph higher up).
#  4. "{{ tag }}" -> "{{ tag|trim(' -F') }}" (curly quotes).
#  5. Turn the stray empty paragraph between the two "{%p endfor %}" lines and
#     "{%p for clause in mysubdoc %}" into a real page-break paragraph.
#  6. Move the "{{r clause[2] }}" / separator / "{{p clause[3] }}" / separator
#     block up so it immediately follows the "{{ clause[0].name }}" heading,
#     followed by a page break and a paragraph holding a new _GoBack bookmark.
#  7. Remove the trailing "Retrieved and assembled at {{ today_date }}."
#     paragraph (and the blank spacer paragraph before it).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "In Word 2019, you can do this by opening the" -> "In Word, you can..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    ". In Word 2019, you can do this by opening the",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ". In Word, you can do this by opening the",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Remove the _GoBack bookmark that sits after "...working on."
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
if ($goBack.Start -ne $goBack.End -or $goBack.Start -ne 0) {
    $goBack.Delete()
}

# ---------------------------------------------------------------------------
# 3. Remove the "You can find an appendix ..." paragraph entirely.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("You can find an appendix of the documents based on the tags that you searched by at the end of the document.")
if ($found) {
    $para = $rng.Paragraphs(1)
    $para.Range.Delete()
}

# ---------------------------------------------------------------------------
# 4. "{{ tag }}:" -> "{{ tag|trim(' -F') }}:" with curly quotes.
# ---------------------------------------------------------------------------
$lsquo = [char]0x2018
$rsquo = [char]0x2019
$tagReplacement = "{{ tag|trim(" + $lsquo + " -F" + $rsquo + ") }}:"
$d.Content.Find.Execute(
    "{{ tag }}:",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    $tagReplacement,
    2) | Out-Null

Write-Host "done-part1"
